$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D15").Value = "Default:192.168.103.34"
$ws.Range("D14").Value = "Default:192.168.103.99"
$ws.Range("D13").Value = "Default:192.168.103.67"

$ws.Range("D13").Select()
